# Fruta / hortaliza, semanal
# Insert 6 new weekly price rows for "Nectarín" (Macroferia Regional de Talca)
# above the existing last block of rows (old rows 816-823 shift down to 822-829),
# then populate the newly opened rows 816-821 with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 816:823 down by 6 rows to make room for the new data.
$ws.Rows.Item(816).Resize(6).Insert()

# Shared/static column values for every row in this data block.
$mercadoId  = 5
$mercado    = "Macroferia Regional de Talca"
$region     = "Maule"
$codreg     = 7
$tipo       = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria  = "Nectarín"
$origen     = "Región de O'Higgins"
$unidad18   = "$/bandeja 18 kilos granel"
$fecha      = 44939

# New rows of weekly data (variety / quality / volume / min / max / avg / $kg / kg per unit).
$newRows = @(
    @{ Row=816; K="Magique";     L="Especial";                M=250; N=17000; O=17000; P=17000; S=944;  T=18 },
    @{ Row=817; K="Magique";     L="Extra (doble especial)";  M=190; N=19000; O=19000; P=19000; S=1056; T=18 },
    @{ Row=818; K="Magique";     L="Segunda";                 M=230; N=15000; O=15000; P=15000; S=833;  T=18 },
    @{ Row=819; K="Super Queen"; L="Especial";                M=150; N=17000; O=17000; P=17000; S=944;  T=18 },
    @{ Row=820; K="Super Queen"; L="Primera";                 M=280; N=15000; O=15000; P=15000; S=833;  T=18 },
    @{ Row=821; K="Super Queen"; L="Segunda";                 M=250; N=12000; O=12000; P=12000; S=667;  T=18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $unidad18
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
